$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# ---------------------------------------------------------------------------
# Insert 4 new rows right above the current last data row (row 70, "Zombie"),
# which will be pushed down to row 74. New rows 70-73 hold 4 new set/style
# combos: Ninja/Agent, Ninja/Martial Artist, Immortal/Angel, Immortal/Devil.
# ---------------------------------------------------------------------------
$ws.Range("A70:Q73").Insert()

# Row 70: Ninja / Agent
$ws.Range("A70").Value = "Ninja"
$ws.Range("B70").Value = "Agent"
$ws.Range("C70").Value = "MS+2"
$ws.Range("D70").Value = "RS+2"
$ws.Range("E70").Value = "RW+2"
$ws.Range("F70").Value = "RS+2"
$ws.Range("G70").Value = "CW+2"

# Row 71: Ninja / Martial Artist
$ws.Range("A71").Value = "Ninja"
$ws.Range("B71").Value = "Martial Artist"
$ws.Range("C71").Value = "MS+3"
$ws.Range("D71").Value = "PW+0"
$ws.Range("E71").Value = "PW+1"
$ws.Range("F71").Value = "RS+2"
$ws.Range("G71").Value = "CW+1"

# Row 72: Immortal / Angel
$ws.Range("A72").Value = "Immortal"
$ws.Range("B72").Value = "Angel"
$ws.Range("C72").Value = "MS+3"
$ws.Range("D72").Value = "PS+3"
$ws.Range("E72").Value = "PW+1"
$ws.Range("F72").Value = "PS+0"
$ws.Range("G72").Value = "CP+3"

# Row 73: Immortal / Devil
$ws.Range("A73").Value = "Immortal"
$ws.Range("B73").Value = "Devil"
$ws.Range("C73").Value = "PR+2"
$ws.Range("D73").Value = "RS+0"
$ws.Range("E73").Value = "RW+2"
$ws.Range("F73").Value = "PW+0"
$ws.Range("G73").Value = "CW+3"

# Fill in the H:Q analysis formulas for the 4 new rows (same pattern used
# throughout the sheet for every data row).
for ($r = 70; $r -le 73; $r++) {
    $ws.Range("H$r").Formula = "=LEN(`$P$r)-LEN(SUBSTITUTE(`$P$r,`"M`",`"`"))"
    $ws.Range("I$r").Formula = "=LEN(`$P$r)-LEN(SUBSTITUTE(`$P$r,`"S`",`"`"))"
    $ws.Range("J$r").Formula = "=LEN(`$P$r)-LEN(SUBSTITUTE(`$P$r,`"R`",`"`"))"
    $ws.Range("K$r").Formula = "=LEN(`$P$r)-LEN(SUBSTITUTE(`$P$r,`"T`",`"`"))"
    $ws.Range("L$r").Formula = "=LEN(`$P$r)-LEN(SUBSTITUTE(`$P$r,`"W`",`"`"))"
    $ws.Range("M$r").Formula = "=LEN(`$P$r)-LEN(SUBSTITUTE(`$P$r,`"C`",`"`"))"
    $ws.Range("N$r").Formula = "=LEN(`$P$r)-LEN(SUBSTITUTE(`$P$r,`"P`",`"`"))"
    $ws.Range("O$r").Formula = "=COUNTIF(H${r}:N$r, `">0`" )"
    $ws.Range("P$r").Formula = "=_xlfn.CONCAT(C${r}:G$r)"
    $ws.Range("Q$r").Formula = "=MAX(H${r}:N$r)"
}

# Match styling used by every other data row: columns A-G use style 1,
# column Q uses style 5 (the rest are left in the sheet's default style,
# same as the row that used to be last).
$ws.Range("A70:G73").Style = $ws.Range("A69").Style
$ws.Range("Q70:Q73").Style = $ws.Range("Q69").Style

# ---------------------------------------------------------------------------
# Defined name: _FilterDatabase needs to track the new data extent.
# ---------------------------------------------------------------------------
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$Q`$75"

# ---------------------------------------------------------------------------
# Update the view selection to match what the author left selected.
# ---------------------------------------------------------------------------
$ws.Range("H73").Select()
